# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# and swaps two pairs of rows (19/20 and 49/50) per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.571.51'
$ws.Range("D3").Value = '2.984.43'
$ws.Range("E3").Value = '  -6.95%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '544.10'
$ws.Range("E5").Value = '  -5.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.90'
$ws.Range("E6").Value = '  -8.59%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.562'
$ws.Range("E8").Value = '  -5.85%  '
$ws.Range("D9").Value = '2.985.62'
$ws.Range("E9").Value = '  -6.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.112'
$ws.Range("E10").Value = '  -6.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.20'
$ws.Range("E11").Value = '  -8.13%  '
$ws.Range("E12").Value = '  -6.94%  '
$ws.Range("D13").Value = '3.507.64'
$ws.Range("E13").Value = '  -6.77%  '
$ws.Range("E14").Value = '  -3.80%  '
$ws.Range("D15").Value = '61.678.63'
$ws.Range("E15").Value = '  -5.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '23.63'
$ws.Range("E16").Value = '  -7.65%  '
$ws.Range("D17").Value = '2.981.62'
$ws.Range("E17").Value = '  -6.90%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000147'
$ws.Range("E18").Value = '  -6.61%  '
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.11'
$ws.Range("E19").Value = '  -4.17%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '388.00'
$ws.Range("E20").Value = '  -5.76%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.89'
$ws.Range("E21").Value = '  -7.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.61'
$ws.Range("E22").Value = '  -7.99%  '
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("E24").Value = '  -6.91%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.468'
$ws.Range("E25").Value = '  -4.67%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.187'
$ws.Range("E26").Value = '  -7.65%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("E27").Value = '  -0.11%  '
$ws.Range("D28").Value = '0.0₃0940'
$ws.Range("E28").Value = '  -10.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.36'
$ws.Range("E29").Value = '  -5.94%  '
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.72'
$ws.Range("E31").Value = '  -6.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.27'
$ws.Range("E32").Value = '  -5.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '158.86'
$ws.Range("E33").Value = '  +1.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.01'
$ws.Range("E34").Value = '  -6.23%  '
$ws.Range("E35").Value = '  -7.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.07'
$ws.Range("E36").Value = '  -6.41%  '
$ws.Range("E37").Value = '  -6.59%  '
$ws.Range("E38").Value = '  -8.55%  '
$ws.Range("D39").Value = '2.440.51'
$ws.Range("E39").Value = '  -11.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.19'
$ws.Range("E40").Value = '  -5.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.87'
$ws.Range("E41").Value = '  -6.66%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.24'
$ws.Range("E42").Value = '  -8.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.659'
$ws.Range("E43").Value = '  -7.62%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0594'
$ws.Range("E44").Value = '  -6.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  -0.10%  '
$ws.Range("E46").Value = '  -6.69%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.94'
$ws.Range("E47").Value = '  -12.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0957'
$ws.Range("E48").Value = '  -3.51%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.71'
$ws.Range("E49").Value = '  -8.24%  '
$ws.Range("B50").Value = 'WhiteBITCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '10.47'
$ws.Range("E50").Value = '  -0.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '264.89'
$ws.Range("E51").Value = '  -10.52%  '
